$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("41:41").Insert()
$ws.Range("A41:E41").ClearFormats()

$ws.Range("A41").Value = "603078cabf4cab0027323ede"
$ws.Range("B41").Value = "Vilmar"
$ws.Range("C41").Value = 25837398.539000001
$ws.Range("D41").Value = 40
$ws.Range("E41").Value = "Vilmarbernardes"

$ws.Range("D42").Value = 41
$ws.Range("D43").Value = 42
$ws.Range("D44").Value = 43
$ws.Range("D45").Value = 44

[void]$ws.Range("A2:E45").Select()

Write-Host "done"
